$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# --- Strip the old per-row distribution formatting (style 8) from rows 24-31 ---
$ws.Range("B24:G31").ClearFormats()
$ws.Range("C24:F24").Value = ""
$ws.Range("C25:G25").Value = ""
$ws.Range("F26:G26").Value = ""
$ws.Range("C27:G30").Value = ""
$ws.Range("C31:G31").Value = ""

# --- Set the new shared-string labels in the exact order the author retyped
#     them, so the workbook's string table compacts/grows in the same order
#     as the authored edit (old "prev_se_*_noninf/_inf" rows 27-30 get
#     retired, freeing their slots before the new labels are appended). ---
$ws.Range("A31").Value = "tpt_completion_perc"
$ws.Range("A27").Value = "prev_se_subclin_noninf_pearl"
$ws.Range("A28").Value = "prev_se_clin_noninf_pearl"
$ws.Range("A29").Value = "prev_se_subclin_inf_pearl"
$ws.Range("A30").Value = "prev_se_clin_inf_pearl"
$ws.Range("A33").Value = "prev_se_clin_noninf_cxr"
$ws.Range("A34").Value = "prev_se_subclin_inf_cxr"
$ws.Range("A35").Value = "prev_se_clin_inf_cxr"
$ws.Range("A32").Value = "prev_se_subclin_noninf_cxr"

# --- Row 24: prev_se_incipient -> 0.75 (distribution cleared, G label kept) ---
$ws.Range("B24").Value = 0.75

# --- Row 25: prev_se_contained -> 0.75 ---
$ws.Range("B25").Value = 0.75

# --- Row 26: prev_se_cleared -> 0.35, uniform(0.2, 0.5) ---
$ws.Range("B26").Value = 0.35
$ws.Range("C26").Value = "uniform"
$ws.Range("D26").Value = 0.2
$ws.Range("E26").Value = 0.5

# --- Row 27: prev_se_subclin_noninf_pearl -> 0.95 ---
$ws.Range("B27").Value = 0.95

# --- Row 28: prev_se_clin_noninf_pearl -> 0.95 ---
$ws.Range("B28").Value = 0.95

# --- Row 29: prev_se_subclin_inf_pearl -> 0.95 ---
$ws.Range("B29").Value = 0.95

# --- Row 30: prev_se_clin_inf_pearl -> 0.95 ---
$ws.Range("B30").Value = 0.95

# --- Row 31: tpt_completion_perc -> 70 ---
$ws.Range("B31").Value = 70

# --- Row 32: prev_se_subclin_noninf_cxr -> 0.5 ---
$ws.Range("B32").Value = 0.5

# --- Row 33: prev_se_clin_noninf_cxr -> 0.9 ---
$ws.Range("B33").Value = 0.9

# --- Row 34: prev_se_subclin_inf_cxr -> 0.6 ---
$ws.Range("B34").Value = 0.6

# --- Row 35: prev_se_clin_inf_cxr -> 0.9 ---
$ws.Range("B35").Value = 0.9

# --- View state: selection + scroll position ---
$ws.Range("C26").Select()
$ws.Application.ActiveWindow.ScrollRow = 16

# --- Column A grows slightly to fit the new, longer labels ---
$ws.Columns("A").ColumnWidth = 25.5
